# Update results with a new "LightGBM" model row for each service group
# and renumber the XGBoost model's trailing result row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new LightGBM rows (previously blank placeholder rows)
# Service 1 (Prime Verifier) -> row 8
$ws.Range("B8").Value = "LightGBM"
$ws.Range("C8").Value = 0.9677
$ws.Range("D8").Value = 0.8472
$ws.Range("E8").Value = -0.1209
$ws.Range("H8").Value = "LightGBM"
$ws.Range("I8").Value = 0.3265
$ws.Range("J8").Value = 0.4122
$ws.Range("K8").Value = 0.3068

# Service 2 (Echo Service) -> row 14
$ws.Range("B14").Value = "LightGBM"
$ws.Range("C14").Value = 0.9988
$ws.Range("D14").Value = 0.9972
$ws.Range("E14").Value = 0.9885
$ws.Range("H14").Value = "LightGBM"
$ws.Range("I14").Value = 0.4671
$ws.Range("J14").Value = 0.4466
$ws.Range("K14").Value = 0.3547

# Service 3 (Hash Generator) -> row 20
$ws.Range("B20").Value = "LightGBM"
$ws.Range("C20").Value = 0.9871
$ws.Range("D20").Value = 0.9956
$ws.Range("E20").Value = 0.9839
$ws.Range("H20").Value = "LightGBM"
$ws.Range("I20").Value = 0.5615
$ws.Range("J20").Value = 0.5649
$ws.Range("K20").Value = 0.4557

# Service 4 (Random Password Generator) -> row 26
$ws.Range("B26").Value = "LightGBM"
$ws.Range("C26").Value = 0.9571
$ws.Range("D26").Value = 0.9729
$ws.Range("E26").Value = 0.961
$ws.Range("H26").Value = "LightGBM"
$ws.Range("I26").Value = 0.8687
$ws.Range("J26").Value = 0.8094
$ws.Range("K26").Value = 0.9785

# Update selection to reflect the last edited cell (K26)
$ws.Range("K26").Select()
